$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where "Физ-ра" (lesson #25 or #27) becomes "Физ-ра1"
$rowsToFizRa1 = @(48, 50, 112, 114, 173, 174, 234, 293, 295)
foreach ($r in $rowsToFizRa1) {
    $ws.Range("B$r").Value = "Физ-ра1"
}

# Rows where "Физ-ра" (lesson #26) becomes "Физ-ра2"
$rowsToFizRa2 = @(49, 113, 235, 294)
foreach ($r in $rowsToFizRa2) {
    $ws.Range("B$r").Value = "Физ-ра2"
}

# Column B now needs its own width (separate from column C) to fit the
# longer "Физ-ра1"/"Физ-ра2" labels.
$ws.Columns.Item(2).ColumnWidth = 14.4

# Selection moved to H13
[void]$ws.Range("H13").Select()
